$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.865.50"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.562.76"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'205.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'21.76"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").Value = "'0.247"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "'0.0864"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "1.785.31"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "1.561.12"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "26.875.97"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "'7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").Value = "'154.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("E31").Value = "  -3.45%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").Value = "1.397.10"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'0.918"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'5.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.88%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "'63.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "1.699.03"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "'86.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").Value = "'0.0505"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("D50").Value = "0.0₇0984"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'0.0953"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.16%  "
